# Commit: "Correct typo on slide 15"
#
# Slide 15 ("Code Review - What Peer Code Review Can Provide") has a bullet
# that reads "Helpful coding techniques by other's" - the apostrophe is a
# typo (it should be the plural "others", not the possessive "other's").
# Find that paragraph wherever it lives on the slide and fix its text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)

$newText = "Helpful coding techniques by others"
$fixed = $false

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($j = 1; $j -le $paraCount; $j++) {
            $para = $tr.Paragraphs($j)
            if ($para.Text.Contains("Helpful coding techniques by other")) {
                $para.Runs(1).Text = $newText
                $fixed = $true
            }
        }
    }
}

if (-not $fixed) {
    # Fall back to the known location (Content Placeholder 2, last bullet)
    # in case the slide text didn't match the expected typo exactly.
    $shape = $s.Shapes.Item("Content Placeholder 2")
    $tr = $shape.TextFrame.TextRange
    $para = $tr.Paragraphs(11)
    $para.Runs(1).Text = $newText
}
